$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the FilesTab Cypher query text in B4 ---
# The "File Type" and "Breed" coalesce lines were dropped from the RETURN
# clause (and the "Association" / "Diagnosis" lines re-indented) as part of
# the ICDC Breed 1-14 script corrections.
$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Samoyed']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
           coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newQuery

# --- Row 4 shrinks now that two lines were removed from the wrapped text ---
$ws.Rows.Item(4).RowHeight = 217.5

# --- Update the saved view state: scrolled down, selection moved to B4 ---
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
